$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.484.93'
$ws.Range('E2').Value = '  -4.69%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.268.01'
$ws.Range('E3').Value = '  -6.14%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.58'
$ws.Range('E5').Value = '  -3.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.62'
$ws.Range('E6').Value = '  -4.94%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').Value = '  -2.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.260.88'
$ws.Range('E9').Value = '  -5.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.185'
$ws.Range('E10').Value = '  -7.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.584'
$ws.Range('E11').Value = '  -4.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.29'
$ws.Range('E12').Value = '  -8.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000263'
$ws.Range('E13').Value = '  -7.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '634.98'
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.53'
$ws.Range('E15').Value = '  -6.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.807.98'
$ws.Range('E16').Value = '  -5.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.575.30'
$ws.Range('E17').Value = '  -4.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.79'
$ws.Range('E18').Value = '  -1.50%  '
$ws.Range('E19').Value = '  -3.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.280.72'
$ws.Range('E20').Value = '  -6.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.32'
$ws.Range('E21').Value = '  -8.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.902'
$ws.Range('E22').Value = '  -4.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.65'
$ws.Range('E23').Value = '  -0.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '106.57'
$ws.Range('E24').Value = '  +7.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.97'
$ws.Range('E25').Value = '  -7.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.97'
$ws.Range('E26').Value = '  -7.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.67'
$ws.Range('E27').Value = '  -6.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.44'
$ws.Range('E28').Value = '  -5.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.67'
$ws.Range('E29').Value = '  -5.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.17'
$ws.Range('E30').Value = '  -6.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.97'
$ws.Range('E31').Value = '  -2.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.30'
$ws.Range('E32').Value = '  -5.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.03'
$ws.Range('E33').Value = '  -4.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '556.26'
$ws.Range('E34').Value = '  +11.53%  '
$ws.Range('E35').Value = '  -3.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '57.29'
$ws.Range('E36').Value = '  -5.56%  '
$ws.Range('B37').Value = 'CoreDAO'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.88'
$ws.Range('E37').Value = '  +42.95%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.612.69'
$ws.Range('E39').Value = '  -2.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.48'
$ws.Range('E40').Value = '  -1.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.70'
$ws.Range('E41').Value = '  -7.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0₃0705'
$ws.Range('E42').Value = '  -10.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.126'
$ws.Range('E43').Value = '  -4.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.340'
$ws.Range('E44').Value = '  -7.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '31.88'
$ws.Range('E45').Value = '  -6.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0414'
$ws.Range('E46').Value = '  -5.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.25'
$ws.Range('E47').Value = '  -2.93%  '
$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.60'
$ws.Range('E48').Value = '  -6.69%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.129'
$ws.Range('E49').Value = '  -3.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.998'
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.25'
$ws.Range('E51').Value = '  +1.83%  '
